# Add hostname-prefixed logging feature
#
# The scan that produced this report re-ran, so every row's "ScanDate"
# (column J) needs to be refreshed to the new run's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

$newScanDate = "2025-09-10 16:44:11"

# Refresh the ScanDate column for every data row that has one.
$ws.Range("J2").Value = $newScanDate
$ws.Range("J4").Value = $newScanDate
$ws.Range("J5").Value = $newScanDate
$ws.Range("J6").Value = $newScanDate
$ws.Range("J7").Value = $newScanDate
$ws.Range("J8").Value = $newScanDate

# Restore the sheet's visible selection.
$ws.Activate()
$ws.Range("E2:L8").Select()
